# TODO.xlsx: remove the "Ecritures de régularisation" TODO-list item from the
# "TODO" sheet (row 21) and shift subsequent rows up, then reapply the
# highlight (green fill) / wrap-text formatting that Excel carries along
# with the remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO")

# Delete the entire row that holds "Ecritures de régularisation" (B21).
# This shifts every row below it up by one, which is exactly what the
# diff shows (old row 22 -> new row 21, ..., old row 25 removed).
$ws.Rows("21").Delete()

# "Modifier le nom de l'état REA en APUR" (now B16) and "Planifier la mise
# en oeuvre de l'application pour Douala" (now B21, formerly B22) pick up
# the same green-fill highlight style used throughout the sheet (e.g. B1).
$ws.Range("B1").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)

# "Dans l'apurement, ajouter une liste déroulante..." (now B22, formerly
# B23) picks up the green-fill + wrap-text style already used on B17.
$ws.Range("B17").Copy()
$ws.Range("B22").PasteSpecial(-4122)

$excel.CutCopyMode = 0
